$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for price cells whose new values would otherwise be parsed as numbers
$textRows = 4,5,6,7,8,9,10,11,15,16,17,18,22,25,26,27,29,31,37,40,41,42,44,45,46,49,50
foreach ($r in $textRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = '28.534.71'
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("D3").Value = '1.560.89'
$ws.Range("E3").Value = '  -0.75%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.28%  '
$ws.Range("D5").Value = '210.24'
$ws.Range("E5").Value = '  -0.80%  '
$ws.Range("D6").Value = '0.487'
$ws.Range("E6").Value = '  -1.24%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.32%  '
$ws.Range("D8").Value = '24.85'
$ws.Range("E8").Value = '  +5.05%  '
$ws.Range("D9").Value = '0.244'
$ws.Range("E9").Value = '  -1.04%  '
$ws.Range("D10").Value = '0.0585'
$ws.Range("E10").Value = '  -0.26%  '
$ws.Range("D11").Value = '0.0896'
$ws.Range("E11").Value = '  +0.29%  '
$ws.Range("D12").Value = '1.784.73'
$ws.Range("E12").Value = '  -0.66%  '
$ws.Range("D13").Value = '1.580.24'
$ws.Range("E13").Value = '  +0.47%  '
$ws.Range("D14").Value = '28.537.80'
$ws.Range("E14").Value = '  +0.55%  '
$ws.Range("D15").Value = '0.513'
$ws.Range("E15").Value = '  -0.33%  '
$ws.Range("D16").Value = '3.63'
$ws.Range("E16").Value = '  -1.36%  '
$ws.Range("D17").Value = '61.21'
$ws.Range("E17").Value = '  -0.62%  '
$ws.Range("D18").Value = '229.40'
$ws.Range("E18").Value = '  +0.48%  '
$ws.Range("E19").Value = '  -0.71%  '
$ws.Range("D20").Value = '0.0₃0677'
$ws.Range("E20").Value = '  -1.04%  '
$ws.Range("E21").Value = '  -0.18%  '
$ws.Range("D22").Value = '3.91'
$ws.Range("E22").Value = '  -0.86%  '
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("E24").Value = '  +1.30%  '
$ws.Range("D25").Value = '150.33'
$ws.Range("E25").Value = '  -0.54%  '
$ws.Range("D26").Value = '14.76'
$ws.Range("E26").Value = '  -1.00%  '
$ws.Range("D27").Value = '0.104'
$ws.Range("E27").Value = '  -0.16%  '
$ws.Range("E28").Value = '  -0.29%  '
$ws.Range("D29").Value = '6.21'
$ws.Range("E29").Value = '  -2.34%  '
$ws.Range("E30").Value = '  -4.29%  '
$ws.Range("D31").Value = '1.05'
$ws.Range("E31").Value = '  -2.63%  '
$ws.Range("E32").Value = '  -0.60%  '
$ws.Range("D33").Value = '1.386.70'
$ws.Range("E33").Value = '  +0.51%  '
$ws.Range("E34").Value = '  -4.43%  '
$ws.Range("E35").Value = '  -3.01%  '
$ws.Range("E36").Value = '  -1.90%  '
$ws.Range("D37").Value = '2.69'
$ws.Range("E37").Value = '  +1.64%  '
$ws.Range("E38").Value = '  -2.64%  '
$ws.Range("E39").Value = '  -1.42%  '
$ws.Range("D40").Value = '1.95'
$ws.Range("E40").Value = '  +2.30%  '
$ws.Range("D41").Value = '0.517'
$ws.Range("E41").Value = '  -0.51%  '
$ws.Range("D42").Value = '0.999'
$ws.Range("E42").Value = '  -0.24%  '
$ws.Range("E43").Value = '  -2.01%  '
$ws.Range("D44").Value = '0.0460'
$ws.Range("E44").Value = '  -2.08%  '
$ws.Range("D45").Value = '63.72'
$ws.Range("E45").Value = '  +2.43%  '
$ws.Range("D46").Value = '5.22'
$ws.Range("E46").Value = '  -2.47%  '
$ws.Range("D47").Value = '1.696.73'
$ws.Range("E47").Value = '  -0.62%  '
$ws.Range("E48").Value = '  -5.37%  '
$ws.Range("D49").Value = '84.97'
$ws.Range("E49").Value = '  -0.47%  '
$ws.Range("D50").Value = '43.21'
$ws.Range("E50").Value = '  +7.41%  '
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '0.0₆0100'
$ws.Range("E51").Value = '  -0.30%  '
